$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 43, pushing the existing rows 43:53 down to 44:54.
$ws.Rows("43:43").Insert()

# Populate the newly inserted row 43 with the new Papaya price record
# (same market/product/category as the surrounding rows, new date/volume/price).
$ws.Cells.Item(43,1).Value = 10
$ws.Cells.Item(43,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(43,3).Value = "La Araucanía"
$ws.Cells.Item(43,4).Value = 44460
$ws.Cells.Item(43,5).Value = 9
$ws.Cells.Item(43,6).Value = "Fruta"
$ws.Cells.Item(43,7).Value = 100108
$ws.Cells.Item(43,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(43,9).Value = 100108004
$ws.Cells.Item(43,10).Value = "Papaya"
$ws.Cells.Item(43,11).Value = "Cultivar IV Región"
$ws.Cells.Item(43,12).Value = "Primera"
$ws.Cells.Item(43,13).Value = 30
$ws.Cells.Item(43,14).Value = 20000
$ws.Cells.Item(43,15).Value = 20000
$ws.Cells.Item(43,16).Value = 20000
$ws.Cells.Item(43,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(43,18).Value = "Provincia del Elquí"
$ws.Cells.Item(43,19).Value = 2000
$ws.Cells.Item(43,20).Value = 10
